# Atualização de bases das ligas: swap the data between pairs of adjacent
# match records that were stored in the wrong order. For each pair of rows
# below, every field except the running index (column A) and the shared
# match Date (column D, identical for both rows of a pair) is exchanged
# between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row1, row2) pairs whose B:AB data (minus Date) must be swapped.
$rowPairs = @(
    @(22, 23),
    @(41, 42),
    @(58, 59),
    @(78, 79),
    @(161, 162)
)

# Columns B .. AB, skipping D (4) which holds the shared match Date and is
# identical for both rows in every pair, so it never needs to move.
$colNums = @()
for ($c = 2; $c -le 28; $c++) {
    if ($c -ne 4) {
        $colNums += $c
    }
}

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($c in $colNums) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        if ($v1 -ne $v2) {
            $cell1.Value2 = $v2
            $cell2.Value2 = $v1
        }
    }
}
